# Apresentacao G5 - atualizacao 1o semestre 2014
# 1) Handout master: atualiza o campo de data em cache (13/03/2014 -> 24/03/2014)
# 2) Insere um novo slide "Escopo do Sistema" antes do slide "Conclusao"
# 3) Ajusta o slide "Conclusao" (agora deslocado) - numero de slide / rodape

$p = $ppt.ActivePresentation

# --- 1) Handout master date (best effort; cached auto-date field) ---
try {
    $hm = $p.HandoutMaster
    $hmDt = $hm.HeadersFooters.DateAndTime
    $hmDt.Visible = $true
    $hmDt.Value = "24/03/2014"
} catch {
}

# --- 2) Novo slide "Escopo do Sistema" inserido na posicao 4 (antes de "Conclusao") ---
$novo = $p.Slides.Add(4, 2)

$novo.Shapes.Item(1).TextFrame.TextRange.Text = "Escopo do Sistema`t"
$novo.Shapes.Item(2).TextFrame.TextRange.Text = "Apresente o escopo do sistema"

$novo.HeadersFooters.DateAndTime.Visible = $true
$novo.HeadersFooters.DateAndTime.UseFormat = $false

$novo.HeadersFooters.Footer.Visible = $true
$novo.HeadersFooters.Footer.Text = "Veris Faculdades TCM/3ADS"

$novo.HeadersFooters.SlideNumber.Visible = $true

$novo.Shapes.Item(3).Name = "Espaço Reservado para Data 3"
$novo.Shapes.Item(4).Name = "Espaço Reservado para Rodapé 4"
$novo.Shapes.Item(5).Name = "Espaço Reservado para Número de Slide 5"

# --- 3) Slide "Conclusao" (agora na posicao 5) - mescla os runs do rodape ---
$concl = $p.Slides.Item(5)
$ftr = $concl.Shapes.Item(5)
$ftr.TextFrame.TextRange.Text = "Metrocamp Faculdades TCM/3ADS"
